# [BAEX-17199]: style css 11
#
# Re-key the 5 driver phone-number cells (B2:B6) from plain numbers to
# text-formatted shared strings (so leading/format digits are preserved as
# entered), and move the active selection to G10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a text number format on B2:B6 first so the subsequent .Value writes
# are stored as shared-string text (t="s") rather than re-parsed as numbers
# -- this is what produces the new numFmtId="49" text style in cellXfs.
$ws.Range("B2:B6").NumberFormat = "@"

$ws.Range("B2").Value = "8425121241"
$ws.Range("B3").Value = "8435121242"
$ws.Range("B4").Value = "8445121243"
$ws.Range("B5").Value = "8455121244"
$ws.Range("B6").Value = "8465121245"

# Move the saved selection from I9 to G10.
$ws.Range("G10").Select()
